# Auto-generated edit script: update TPM values and cluster labels (Fgf15-Fgfr2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Fgf15"
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.018508
$ws.Range("H2").Value = 0.055524
$ws.Range("I2").Value = 0.2347866901774728
$ws.Range("J2").Value = 0.3151800006811757
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2858606666666667
$ws.Range("N2").Value = 0.857582
$ws.Range("O2").Value = 0.0687156860066334
$ws.Range("P2").Value = 0.06932858672617494
$ws.Range("Q2").Value = 0.005290709218666667
$ws.Range("R2").Value = 0.047616382968
$ws.Range("S2").Value = 0.01613352848077194
$ws.Range("T2").Value = 0.02185098401158076

# Row 3
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Fgf15"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.018508
$ws.Range("H3").Value = 0.055524
$ws.Range("I3").Value = 0.2347866901774728
$ws.Range("J3").Value = 0.3151800006811757
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.763360333333333
$ws.Range("N3").Value = 11.290081
$ws.Range("O3").Value = 0.90464312565499
$ws.Range("P3").Value = 0.9127119736118995
$ws.Range("Q3").Value = 0.06965227304933333
$ws.Range("R3").Value = 0.626870457444
$ws.Range("S3").Value = 0.2123981652643387
$ws.Range("T3").Value = 0.2876685604647157

# Row 4
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Fgf15"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.018508
$ws.Range("H4").Value = 0.055524
$ws.Range("I4").Value = 0.2347866901774728
$ws.Range("J4").Value = 0.3151800006811757
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.110331
$ws.Range("N4").Value = 0.220662
$ws.Range("O4").Value = 0.02652155835639462
$ws.Range("P4").Value = 0.01783874265571248
$ws.Range("Q4").Value = 0.002042006148
$ws.Range("R4").Value = 0.012252036888
$ws.Range("S4").Value = 0.006226908904846587
$ws.Range("T4").Value = 0.005622414922378776

# Row 5
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Fgf15"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.018508
$ws.Range("H5").Value = 0.055524
$ws.Range("I5").Value = 0.2347866901774728
$ws.Range("J5").Value = 0.3151800006811757
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.0004976666666666667
$ws.Range("N5").Value = 0.001493
$ws.Range("O5").Value = 0.0001196299819817856
$ws.Range("P5").Value = 0.0001206970062130259
$ws.Range("Q5").Value = 0.000009210814666666666
$ws.Range("R5").Value = 0.00008289733199999999
$ws.Range("S5").Value = 0.00002808752751549415
$ws.Range("T5").Value = 0.00003804128250043737

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Fgf15"
$ws.Range("C6").Value = "Fgfr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.060321
$ws.Range("H6").Value = 0.120642
$ws.Range("I6").Value = 0.7652133098225272
$ws.Range("J6").Value = 0.6848199993188243
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2858606666666667
$ws.Range("N6").Value = 0.857582
$ws.Range("O6").Value = 0.0687156860066334
$ws.Range("P6").Value = 0.06932858672617494
$ws.Range("Q6").Value = 0.017243401274
$ws.Range("R6").Value = 0.103460407644
$ws.Range("S6").Value = 0.05258215752586146
$ws.Range("T6").Value = 0.04747760271459418

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Fgf15"
$ws.Range("C7").Value = "Fgfr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.060321
$ws.Range("H7").Value = 0.120642
$ws.Range("I7").Value = 0.7652133098225272
$ws.Range("J7").Value = 0.6848199993188243
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.763360333333333
$ws.Range("N7").Value = 11.290081
$ws.Range("O7").Value = 0.90464312565499
$ws.Range("P7").Value = 0.9127119736118995
$ws.Range("Q7").Value = 0.227009658667
$ws.Range("R7").Value = 1.362057952002
$ws.Range("S7").Value = 0.6922449603906513
$ws.Range("T7").Value = 0.6250434131471838

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Fgf15"
$ws.Range("C8").Value = "Fgfr2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.060321
$ws.Range("H8").Value = 0.120642
$ws.Range("I8").Value = 0.7652133098225272
$ws.Range("J8").Value = 0.6848199993188243
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.110331
$ws.Range("N8").Value = 0.220662
$ws.Range("O8").Value = 0.02652155835639462
$ws.Range("P8").Value = 0.01783874265571248
$ws.Range("Q8").Value = 0.006655276251
$ws.Range("R8").Value = 0.026621105004
$ws.Range("S8").Value = 0.02029464945154803
$ws.Range("T8").Value = 0.0122163277333337

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Fgf15"
$ws.Range("C9").Value = "Fgfr2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.060321
$ws.Range("H9").Value = 0.120642
$ws.Range("I9").Value = 0.7652133098225272
$ws.Range("J9").Value = 0.6848199993188243
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.0004976666666666667
$ws.Range("N9").Value = 0.001493
$ws.Range("O9").Value = 0.0001196299819817856
$ws.Range("P9").Value = 0.0001206970062130259
$ws.Range("Q9").Value = 0.000030019751
$ws.Range("R9").Value = 0.000180118506
$ws.Range("S9").Value = 0.00009154245446629145
$ws.Range("T9").Value = 0.00008265572371258853

